# Update cryptos list values per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "58.004.56"
Set-TextValue "E2" "  -1.82%  "
Set-TextValue "D3" "2.469.29"
Set-TextValue "E3" "  -2.24%  "
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "518.43"
Set-TextValue "E5" "  -3.51%  "
Set-TextValue "D6" "131.05"
Set-TextValue "E6" "  -4.35%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  +0.11%  "
Set-TextValue "E8" "  -2.30%  "
Set-TextValue "D9" "0.0991"
Set-TextValue "E9" "  -2.20%  "
Set-TextValue "D10" "0.157"
Set-TextValue "E10" "  -0.55%  "
Set-TextValue "D11" "5.34"
Set-TextValue "E12" "  -1.44%  "
Set-TextValue "D13" "2.905.59"
Set-TextValue "E13" "  -1.30%  "
Set-TextValue "D14" "57.919.26"
Set-TextValue "E14" "  -1.70%  "
Set-TextValue "D15" "22.29"
Set-TextValue "E15" "  -3.34%  "
Set-TextValue "E16" "  -2.29%  "
Set-TextValue "D17" "2.475.60"
Set-TextValue "E17" "  -1.35%  "
Set-TextValue "D18" "10.80"
Set-TextValue "E18" "  -3.17%  "
Set-TextValue "E19" "  -2.62%  "
Set-TextValue "D20" "318.71"
Set-TextValue "E20" "  -1.68%  "
Set-TextValue "E21" "  +0.02%  "
Set-TextValue "D22" "5.74"
Set-TextValue "E22" "  -3.64%  "
Set-TextValue "D23" "64.08"
Set-TextValue "E23" "  -2.78%  "
Set-TextValue "D24" "0.409"
Set-TextValue "E24" "  -3.22%  "
Set-TextValue "E25" "  +0.38%  "
Set-TextValue "E26" "  -3.33%  "
Set-TextValue "D27" "7.33"
Set-TextValue "E27" "  -2.99%  "
Set-TextValue "D28" "0.0₃0750"
Set-TextValue "E28" "  -2.96%  "
Set-TextValue "B29" "Aptos"
Set-TextValue "C29" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D29" "6.31"
Set-TextValue "E29" "  -5.86%  "
Set-TextValue "B30" "Monero"
Set-TextValue "C30" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D30" "166.01"
Set-TextValue "E30" "  -1.27%  "
Set-TextValue "E31" "  -4.54%  "
Set-TextValue "E32" "  -2.15%  "
Set-TextValue "E33" "  +0.11%  "
Set-TextValue "E34" "  -0.01%  "
Set-TextValue "D35" "18.06"
Set-TextValue "E35" "  -2.07%  "
Set-TextValue "E36" "  -10.47%  "
Set-TextValue "E37" "  -3.38%  "
Set-TextValue "E38" "  -4.70%  "
Set-TextValue "D39" "0.791"
Set-TextValue "E39" "  -2.72%  "
Set-TextValue "E40" "  -4.45%  "
Set-TextValue "D41" "272.25"
Set-TextValue "E41" "  -4.35%  "
Set-TextValue "D42" "5.01"
Set-TextValue "E42" "  -2.78%  "
Set-TextValue "E43" "  -2.57%  "
Set-TextValue "D44" "126.42"
Set-TextValue "E44" "  -4.87%  "
Set-TextValue "E45" "  -2.42%  "
Set-TextValue "E46" "  -4.03%  "
Set-TextValue "E47" "  -3.25%  "
Set-TextValue "D48" "17.08"
Set-TextValue "E48" "  -1.69%  "
Set-TextValue "D49" "1.731.65"
Set-TextValue "E49" "  -1.98%  "
Set-TextValue "D50" "0.975"
Set-TextValue "E50" "  -1.28%  "
Set-TextValue "E51" "  -1.15%  "
